$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.98"
$ws.Range("D2").Style = $ws.Range("D2").Style
$ws.Range("E2").Value = "'0.25%"
$ws.Range("E2").Style = $ws.Range("E2").Style
$ws.Range("D3").Value = "'41.10"
$ws.Range("D3").Style = $ws.Range("D3").Style
$ws.Range("E3").Value = "'1.62%"
$ws.Range("E3").Style = $ws.Range("E3").Style
$ws.Range("E4").Value = "'-1.87%"
$ws.Range("E4").Style = $ws.Range("E4").Style
$ws.Range("D5").Value = "'0.08066"
$ws.Range("D5").Style = $ws.Range("D5").Style
$ws.Range("D6").Value = "'2.020"
$ws.Range("D6").Style = $ws.Range("D6").Style
$ws.Range("E6").Value = "'2.98%"
$ws.Range("E6").Style = $ws.Range("E6").Style
$ws.Range("D7").Value = "'8.722"
$ws.Range("D7").Style = $ws.Range("D7").Style
$ws.Range("E7").Value = "'-0.66%"
$ws.Range("E7").Style = $ws.Range("E7").Style
$ws.Range("D8").Value = "'4.512"
$ws.Range("D8").Style = $ws.Range("D8").Style
$ws.Range("E8").Value = "'-1.83%"
$ws.Range("E8").Style = $ws.Range("E8").Style
$ws.Range("D10").Value = "'0.9223"
$ws.Range("D10").Style = $ws.Range("D10").Style
$ws.Range("E10").Value = "'-2.32%"
$ws.Range("E10").Style = $ws.Range("E10").Style
$ws.Range("D11").Value = "'0.1273"
$ws.Range("D11").Style = $ws.Range("D11").Style
$ws.Range("E11").Value = "'-0.79%"
$ws.Range("E11").Style = $ws.Range("E11").Style
$ws.Range("D12").Value = "'0.1943"
$ws.Range("D12").Style = $ws.Range("D12").Style
$ws.Range("E12").Value = "'-2.16%"
$ws.Range("E12").Style = $ws.Range("E12").Style
$ws.Range("E13").Value = "'-8.20%"
$ws.Range("E13").Style = $ws.Range("E13").Style
$ws.Range("D14").Value = "'0.09374"
$ws.Range("D14").Style = $ws.Range("D14").Style
$ws.Range("E14").Value = "'1.79%"
$ws.Range("E14").Style = $ws.Range("E14").Style
$ws.Range("D15").Value = "'0.03701"
$ws.Range("D15").Style = $ws.Range("D15").Style
$ws.Range("E15").Value = "'5.02%"
$ws.Range("E15").Style = $ws.Range("E15").Style
$ws.Range("D16").Value = "'0.1052"
$ws.Range("D16").Style = $ws.Range("D16").Style
$ws.Range("E16").Value = "'9.32%"
$ws.Range("E16").Style = $ws.Range("E16").Style
$ws.Range("D17").Value = "'0.001295"
$ws.Range("D17").Style = $ws.Range("D17").Style
$ws.Range("E17").Value = "'-2.48%"
$ws.Range("E17").Style = $ws.Range("E17").Style
$ws.Range("D18").Value = "'0.006239"
$ws.Range("D18").Style = $ws.Range("D18").Style
$ws.Range("E18").Value = "'-3.98%"
$ws.Range("E18").Style = $ws.Range("E18").Style
$ws.Range("D19").Value = "'3.364"
$ws.Range("D19").Style = $ws.Range("D19").Style
$ws.Range("E19").Value = "'-0.17%"
$ws.Range("E19").Style = $ws.Range("E19").Style
$ws.Range("E21").Value = "'-1.08%"
$ws.Range("E21").Style = $ws.Range("E21").Style
$ws.Range("D22").Value = "'0.2651"
$ws.Range("D22").Style = $ws.Range("D22").Style
$ws.Range("E22").Value = "'9.80%"
$ws.Range("E22").Style = $ws.Range("E22").Style
$ws.Range("D23").Value = "'0.04428"
$ws.Range("D23").Style = $ws.Range("D23").Style
$ws.Range("E23").Value = "'0.12%"
$ws.Range("E23").Style = $ws.Range("E23").Style
$ws.Range("D24").Value = "'0.001262"
$ws.Range("D24").Style = $ws.Range("D24").Style
$ws.Range("E24").Value = "'0.19%"
$ws.Range("E24").Style = $ws.Range("E24").Style
$ws.Range("D25").Value = "'0.004393"
$ws.Range("D25").Style = $ws.Range("D25").Style
$ws.Range("E25").Value = "'0.60%"
$ws.Range("E25").Style = $ws.Range("E25").Style
$ws.Range("D26").Value = "'0.0001240"
$ws.Range("D26").Style = $ws.Range("D26").Style
$ws.Range("E26").Value = "'8.58%"
$ws.Range("E26").Style = $ws.Range("E26").Style
$ws.Range("D39").Value = "'0.02845"
$ws.Range("D39").Style = $ws.Range("D39").Style
$ws.Range("E39").Value = "'16.40%"
$ws.Range("E39").Style = $ws.Range("E39").Style
$ws.Range("D40").Value = "'0.05463"
$ws.Range("D40").Style = $ws.Range("D40").Style
$ws.Range("E40").Value = "'3.03%"
$ws.Range("E40").Style = $ws.Range("E40").Style
$ws.Range("D41").Value = "'0.007622"
$ws.Range("D41").Style = $ws.Range("D41").Style
$ws.Range("E41").Value = "'2.14%"
$ws.Range("E41").Style = $ws.Range("E41").Style
$ws.Range("D42").Value = "'0.009947"
$ws.Range("D42").Style = $ws.Range("D42").Style
$ws.Range("E42").Value = "'14.44%"
$ws.Range("E42").Style = $ws.Range("E42").Style
$ws.Range("D43").Value = "'0.1419"
$ws.Range("D43").Style = $ws.Range("D43").Style
$ws.Range("E43").Value = "'-0.87%"
$ws.Range("E43").Style = $ws.Range("E43").Style
$ws.Range("D44").Value = "'0.002130"
$ws.Range("D44").Style = $ws.Range("D44").Style
$ws.Range("E44").Value = "'0.29%"
$ws.Range("E44").Style = $ws.Range("E44").Style
$ws.Range("D45").Value = "'0.01188"
$ws.Range("D45").Style = $ws.Range("D45").Style
$ws.Range("E45").Value = "'9.91%"
$ws.Range("E45").Style = $ws.Range("E45").Style
$ws.Range("D46").Value = "'0.00006768"
$ws.Range("D46").Style = $ws.Range("D46").Style
$ws.Range("E46").Value = "'-1.55%"
$ws.Range("E46").Style = $ws.Range("E46").Style
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = $ws.Range("D47").Style
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = $ws.Range("E47").Style
$ws.Range("D48").Value = "'0.002999"
$ws.Range("D48").Style = $ws.Range("D48").Style
$ws.Range("E48").Value = "'-5.55%"
$ws.Range("E48").Style = $ws.Range("E48").Style
$ws.Range("D49").Value = "'0.002280"
$ws.Range("D49").Style = $ws.Range("D49").Style
$ws.Range("E49").Value = "'33.83%"
$ws.Range("E49").Style = $ws.Range("E49").Style
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = $ws.Range("D50").Style
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = $ws.Range("E50").Style
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = $ws.Range("D51").Style
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = $ws.Range("E51").Style
